$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: add date text in A13, quantity 8 in D13, formula in F13
$ws.Range("A13").Value = "26.08.2022"
$ws.Range("D13").Value = 8
$ws.Range("F13").Formula = "=B13-D13"

# Row 14: add formula referencing F13
$ws.Range("B14").Formula = "=F13"

# Restore row 11 height to default (remove explicit row height override)
$ws.Rows.Item(11).AutoFit()

$ws.Range("B13:C14").Select()
